# Updates the crypto price/volume table (columns D and E) to the latest scrape.
# Values in column D ("Price") can look numeric (e.g. "275.90"); assigning them
# plainly through COM would let Excel coerce them to real numbers (dropping
# trailing zeros / switching the stored cell type). Forcing NumberFormat "@"
# (Text) before the write, then restoring the "Normal" style afterwards, keeps
# every cell a literal text value without leaving any stray formatting behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '51.642.56'
Set-TextValue "E2" '  +3.65%  '
Set-TextValue "D3" '2.738.93'
Set-TextValue "E3" '  +2.46%  '
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  -0.03%  '
Set-TextValue "D5" '115.47'
Set-TextValue "E5" '  +1.78%  '
Set-TextValue "D6" '333.13'
Set-TextValue "E6" '  +2.26%  '
Set-TextValue "E7" '  +0.40%  '
Set-TextValue "E8" '  -0.02%  '
Set-TextValue "E9" '  +4.14%  '
Set-TextValue "E10" '  +1.68%  '
Set-TextValue "D11" '20.15'
Set-TextValue "E11" '  +0.24%  '
Set-TextValue "D12" '0.0827'
Set-TextValue "E12" '  +0.63%  '
Set-TextValue "E13" '  +2.54%  '
Set-TextValue "E14" '  +3.28%  '
Set-TextValue "D15" '3.169.70'
Set-TextValue "E15" '  +2.50%  '
Set-TextValue "D16" '2.735.39'
Set-TextValue "E16" '  +2.80%  '
Set-TextValue "E17" '  +1.24%  '
Set-TextValue "D18" '51.528.13'
Set-TextValue "E18" '  +3.46%  '
Set-TextValue "D19" '13.82'
Set-TextValue "E19" '  +5.34%  '
Set-TextValue "E20" '  +1.14%  '
Set-TextValue "E21" '  +1.83%  '
Set-TextValue "D22" '0.0₃0959'
Set-TextValue "E22" '  -0.15%  '
Set-TextValue "D23" '275.90'
Set-TextValue "E23" '  -0.37%  '
Set-TextValue "D24" '69.95'
Set-TextValue "E24" '  -2.48%  '
Set-TextValue "E25" '  +3.85%  '
Set-TextValue "D26" '26.81'
Set-TextValue "E26" '  +0.01%  '
Set-TextValue "E27" '  +0.60%  '
Set-TextValue "E28" '  +0.10%  '
Set-TextValue "D29" '10.32'
Set-TextValue "E29" '  +1.41%  '
Set-TextValue "E30" '  -0.84%  '
Set-TextValue "D31" '35.71'
Set-TextValue "E31" '  -1.28%  '
Set-TextValue "E32" '  -1.57%  '
Set-TextValue "D33" '50.32'
Set-TextValue "E33" '  +0.15%  '
Set-TextValue "E34" '  +2.22%  '
Set-TextValue "D35" '0.0823'
Set-TextValue "E35" '  +1.91%  '
Set-TextValue "D36" '19.37'
Set-TextValue "E36" '  -0.77%  '
Set-TextValue "D37" '0.999'
Set-TextValue "E37" '  -0.26%  '
Set-TextValue "E38" '  +1.30%  '
Set-TextValue "D39" '4.99'
Set-TextValue "E39" '  -0.82%  '
Set-TextValue "E40" '  +2.39%  '
Set-TextValue "D41" '23.90'
Set-TextValue "E41" '  +6.15%  '
Set-TextValue "D42" '128.41'
Set-TextValue "E42" '  +2.39%  '
Set-TextValue "D43" '0.0347'
Set-TextValue "E44" '  +3.48%  '
Set-TextValue "E45" '  +0.52%  '
Set-TextValue "D46" '2.38'
Set-TextValue "E46" '  +16.21%  '
Set-TextValue "D47" '2.096.49'
Set-TextValue "E47" '  -1.08%  '
Set-TextValue "E48" '  +1.66%  '
Set-TextValue "E49" '  +1.03%  '
Set-TextValue "D50" '5.59'
Set-TextValue "E50" '  +4.92%  '
Set-TextValue "D51" '8.98'
Set-TextValue "E51" '  -0.61%  '
